# bug of signup is solved
# - Remove the "Users" sheet entirely (it held the broken/duplicate signup
#   test row).
# - Keep "Sheet1" as the only worksheet, lower-case its header row, and
#   refresh its data rows with the corrected signup records.

$wb = $excel.ActiveWorkbook

# Delete the "Users" worksheet.
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Delete()

# Remaining worksheet.
$ws = $wb.Worksheets.Item("Sheet1")

# Lower-case the header row.
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "email"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "department"

# Row 2
$ws.Range("A2").Value = "manvir"
$ws.Range("B2").Value = "g@gmail.com"
$ws.Range("C2").Value = "`$2b`$10`$wsiUzzuFg/VERhNRUvu72uJl.XTEiFH4P16D/cUxeBpYXIDz7ZkE6"
$ws.Range("D2").Value = "forensic"

# Row 3
$ws.Range("A3").Value = "manvir1"
$ws.Range("B3").Value = "ga@gmail.com"
$ws.Range("C3").Value = "`$2b`$10`$LmUktpKzoDl3RXiWEnrxoOJU7zTotTiPqiljQk6gWz88p3e6ab2oG"
$ws.Range("D3").Value = "forensic"

# Row 4
$ws.Range("A4").Value = "mavnr"
$ws.Range("B4").Value = "a@gmail.com"
$ws.Range("C4").Value = "`$2b`$10`$./0i/7lnt2iUMfH3hyRz8.8QOB.Yx/WYVNU5SNXtb709ZlXo8olFC"
$ws.Range("D4").Value = "academics"
